$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row 17 data for "Exp 21"
$ws.Range("A17").Value = "Exp 21"
$ws.Range("B17").Value = 0.35
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = "Local"
$ws.Range("E17").Value = -1
$ws.Range("F17").Value = "Exp 21.png"

# Apply same style (centered alignment) as the other rows to A17:E17
$ws.Range("A17:E17").HorizontalAlignment = -4108

# Match the row height used by the preceding rows (14.25pt, custom height)
$ws.Rows.Item(17).RowHeight = 14.25

# Update the selection as recorded in the saved workbook
$ws.Range("I15").Select()
